$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "36.082.14"
Set-TextValue "E2" "  -3.87%  "

# Row 3
Set-TextValue "D3" "1.963.39"
Set-TextValue "E3" "  -4.37%  "

# Row 4
Set-TextValue "E4" "  -0.20%  "

# Row 5
Set-TextValue "D5" "242.39"
Set-TextValue "E5" "  -3.95%  "

# Row 6
Set-TextValue "D6" "0.627"
Set-TextValue "E6" "  -3.12%  "

# Row 7
Set-TextValue "D7" "62.60"
Set-TextValue "E7" "  -5.71%  "

# Row 8
Set-TextValue "E8" "  -0.18%  "

# Row 9
Set-TextValue "D9" "0.372"
Set-TextValue "E9" "  -1.18%  "

# Row 10
Set-TextValue "D10" "56.09"
Set-TextValue "E10" "  -5.69%  "

# Row 11
Set-TextValue "D11" "0.0807"
Set-TextValue "E11" "  +6.84%  "

# Row 12
Set-TextValue "E12" "  -1.20%  "

# Row 13
Set-TextValue "D13" "0.860"
Set-TextValue "E13" "  -5.79%  "

# Row 14
Set-TextValue "D14" "22.24"
Set-TextValue "E14" "  +7.34%  "

# Row 15
Set-TextValue "D15" "14.07"
Set-TextValue "E15" "  -7.78%  "

# Row 16
Set-TextValue "D16" "2.247.97"
Set-TextValue "E16" "  -4.54%  "

# Row 17
Set-TextValue "D17" "5.43"
Set-TextValue "E17" "  -3.43%  "

# Row 18
Set-TextValue "D18" "1.964.11"
Set-TextValue "E18" "  -4.26%  "

# Row 19
Set-TextValue "D19" "35.955.86"
Set-TextValue "E19" "  -4.09%  "

# Row 20
Set-TextValue "D20" "71.11"
Set-TextValue "E20" "  -3.21%  "

# Row 21
Set-TextValue "D21" "0.0₃0856"
Set-TextValue "E21" "  -2.66%  "

# Row 22
Set-TextValue "D22" "237.66"
Set-TextValue "E22" "  -0.04%  "

# Row 23
Set-TextValue "D23" "5.21"
Set-TextValue "E23" "  -2.78%  "

# Row 24
Set-TextValue "E24" "  +0.28%  "

# Row 25
Set-TextValue "D25" "2.55"
Set-TextValue "E25" "  -7.27%  "

# Row 26
Set-TextValue "D26" "2.30"
Set-TextValue "E26" "  -2.27%  "

# Row 27
Set-TextValue "D27" "9.84"
Set-TextValue "E27" "  +2.65%  "

# Row 28
Set-TextValue "D28" "159.46"
Set-TextValue "E28" "  -3.68%  "

# Row 29
Set-TextValue "D29" "19.85"
Set-TextValue "E29" "  -0.34%  "

# Row 30
Set-TextValue "D30" "0.131"
Set-TextValue "E30" "  +17.85%  "

# Row 31
Set-TextValue "E31" "  -1.66%  "

# Row 32
Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.14"
Set-TextValue "E32" "  -6.34%  "

# Row 33
Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "4.85"
Set-TextValue "E33" "  -7.46%  "

# Row 34
Set-TextValue "D34" "0.0621"
Set-TextValue "E34" "  +0.92%  "

# Row 35
Set-TextValue "D35" "4.41"
Set-TextValue "E35" "  -6.85%  "

# Row 36
Set-TextValue "D36" "6.30"
Set-TextValue "E36" "  +4.46%  "

# Row 37
Set-TextValue "B37" "LidoDAOToken"
Set-TextValue "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D37" "2.29"
Set-TextValue "E37" "  -6.83%  "

# Row 38
Set-TextValue "B38" "BinanceUSD"
Set-TextValue "C38" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D38" "1.00"
Set-TextValue "E38" "  -0.22%  "

# Row 39
Set-TextValue "D39" "1.83"
Set-TextValue "E39" "  +0.85%  "

# Row 40
Set-TextValue "D40" "3.12"
Set-TextValue "E40" "  +15.45%  "

# Row 41
Set-TextValue "D41" "0.0985"
Set-TextValue "E41" "  -4.96%  "

# Row 42
Set-TextValue "E42" "  -1.18%  "

# Row 43
Set-TextValue "D43" "0.0212"
Set-TextValue "E43" "  -3.27%  "

# Row 44
Set-TextValue "E44" "  -4.70%  "

# Row 45
Set-TextValue "D45" "1.09"
Set-TextValue "E45" "  -4.68%  "

# Row 46
Set-TextValue "B46" "InjectiveProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "16.22"
Set-TextValue "E46" "  -5.04%  "

# Row 47
Set-TextValue "B47" "Aave"
Set-TextValue "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "92.50"
Set-TextValue "E47" "  -3.16%  "

# Row 48
Set-TextValue "D48" "7.56"
Set-TextValue "E48" "  -7.12%  "

# Row 49
Set-TextValue "D49" "1.338.78"
Set-TextValue "E49" "  -6.35%  "

# Row 50
Set-TextValue "D50" "2.79"
Set-TextValue "E50" "  -5.47%  "

# Row 51
Set-TextValue "D51" "2.140.37"
Set-TextValue "E51" "  -4.58%  "

